# Home Work2 Final Commit
# Updates the results table: widen the "approxSilhFull" column slightly and
# refresh the measured accuracy / timing values that were re-run.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Widen the 4th grid column (approxSilhFull) from 2168 -> 2279 twips ---
$t.Columns.Item(4).Width = 2279 / 20

# --- Refresh the measured values in the results table ---
# (MatchWholeWord = $true keeps short "NNms" values from matching inside
#  longer ones, e.g. "34ms" inside "15334ms".)

function Replace-Value($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Row: Uber_10_large.csv (k=10, t=2000)
Replace-Value "0.9999999999977554" "0.6284781034778061"
Replace-Value "56565ms"            "140525ms"
Replace-Value "0.943447226766584"  "0.5284811736171426"
Replace-Value "4948ms"             "5142ms"

# Row: Uber_3_large.csv (k=3, t=2000)
Replace-Value "0.9290682217994826" "0.49262204165251927"
Replace-Value "15334ms"            "17583ms"
Replace-Value "0.9263242530755712" "0.4777229036647352"
Replace-Value "466ms"              "631ms"

# Row: Uber_3_small.csv (k=3, t=40)
Replace-Value "0.9310338429411754" "0.8782979000141207"
Replace-Value "59ms"               "80ms"
Replace-Value "0.5114864864864865" "0.808918819124161"
Replace-Value "1ms"                "4ms"

# Row: Uber_3_small.csv (k=3, t=80)
Replace-Value "0.9299437989409424" "0.916279244530017"
Replace-Value "54ms"               "50ms"
Replace-Value "0.9935064935064936" "0.8204577684646315"

# Row: Uber_3_small.csv (k=3, t=1000)
Replace-Value "0.9289446866612711" "0.9148194891306463"
Replace-Value "105ms"              "110ms"
Replace-Value "0.9990118577075099" "0.9148194891306459"
Replace-Value "34ms"               "33ms"
